$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: delete every hyperlink on a worksheet whose current Range address
# matches one of the given addresses (e.g. '$A$3'). Re-queries the live
# collection after every deletion since deleting one hyperlink can make
# stale references to the rest of the collection unreliable.
# ---------------------------------------------------------------------------
function Remove-HyperlinksAt($ws, $addresses) {
    foreach ($target in $addresses) {
        $found = $true
        while ($found) {
            $found = $false
            foreach ($hl in $ws.Hyperlinks) {
                if ($hl.Range.Address() -eq $target) {
                    $hl.Delete()
                    $found = $true
                    break
                }
            }
        }
    }
}

# ===========================================================================
# Sheet "Overview" (sheet1): File Name | zh-cn | de-de
#   Row2 = 84296453-...md   -> status changes "Handed back" -> "Not yet handed off"
#   Row3 = e8b00593-...md   -> removed entirely
#   Row4 = .localization-config -> shifts up to row3
# ===========================================================================
$wsOverview = $wb.Worksheets.Item(1)

# Drop the whole "e8b00593" row; everything below shifts up automatically.
$wsOverview.Rows.Item(3).Delete()

# Update the handoff status text for the remaining tracked file.
$wsOverview.Range("B2").Value = "Not yet handed off"
$wsOverview.Range("C2").Value = "Not yet handed off"

# The hyperlink collection does not follow the row shift automatically, so
# clean up the stale entries (old row3 = e8b00593, old row4 = .localization-config)
# and recreate the one that is still needed, now anchored on row 3.
Remove-HyperlinksAt $wsOverview @('$A$3', '$A$4')
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/3d931fae12125d8d8549e0f5a9c2c9ccb4a63475/.localization-config", "", "", ".localization-config") | Out-Null

# ===========================================================================
# Sheet "zh-cn" (sheet2)
#   Row2 = 84296453 file details -> status text + handoff datetime updated
#   Row3 = e8b00593 file details -> removed entirely
#   Row4 = .localization-config  -> shifts up to row3
# ===========================================================================
$wsZhCn = $wb.Worksheets.Item(2)

$wsZhCn.Rows.Item(3).Delete()

$wsZhCn.Range("B2").Value = "Not yet handed off"
$wsZhCn.Range("D2").Value = "2016-01-11 02:44:19"

Remove-HyperlinksAt $wsZhCn @('$A$3', '$C$3', '$E$3', '$F$3', '$A$4')
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/3d931fae12125d8d8549e0f5a9c2c9ccb4a63475/.localization-config", "", "", ".localization-config") | Out-Null

# ===========================================================================
# Sheet "de-de" (sheet3)
#   Row2 = 84296453 file details -> status text + handoff datetime updated
#   Row3 = e8b00593 file details -> removed entirely
#   Row4 = .localization-config  -> shifts up to row3
# ===========================================================================
$wsDeDe = $wb.Worksheets.Item(3)

$wsDeDe.Rows.Item(3).Delete()

$wsDeDe.Range("B2").Value = "Not yet handed off"
$wsDeDe.Range("D2").Value = "2016-01-11 02:44:33"

Remove-HyperlinksAt $wsDeDe @('$A$3', '$C$3', '$E$3', '$F$3', '$A$4')
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/3d931fae12125d8d8549e0f5a9c2c9ccb4a63475/.localization-config", "", "", ".localization-config") | Out-Null

Write-Output "done"
